# Populate the (previously empty) Sheet1 with the daily driver report
# header row and a sample data row, matching the exported report layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): bold, thin border all around, centered / top aligned ---
$headers = @("name", "employee_id", "asset", "arrival", "status", "division", "job_title")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous (thin)
}

# --- Data row (row 2): plain values, default formatting ---
$data = @("Roger Doddy", "DODROG", "PT-07S", "04:45 AM", "On Time", "TEXDIST", "Select Maintenance Employee")
for ($i = 0; $i -lt $data.Count; $i++) {
    $cell = $ws.Cells.Item(2, $i + 1)
    $cell.Value = $data[$i]
}
